$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 46

# Row 6
$ws.Range("A6").Value = 41

# Row 7
$ws.Range("A7").Value = 56

# Row 9
$ws.Range("A9").Value = 12

# Row 11
$ws.Range("F11").Value = 473

# Row 14
$ws.Range("A14").Value = 48

# Row 17
$ws.Range("A17").Value = 53

# Row 18
$ws.Range("A18").Value = 57

# Row 20
$ws.Range("A20").Value = 45
$ws.Range("C20").Value = "Fanfare Tickets"
$ws.Range("E20").Value = 4.4
$ws.Range("F20").Value = 14

# Row 21
$ws.Range("A21").Value = 24
$ws.Range("C21").Value = "Fern"
$ws.Range("E21").Value = 4.3
$ws.Range("F21").Value = 34

# Row 22
$ws.Range("A22").Value = 43
$ws.Range("C22").Value = "Floor & Decor"
$ws.Range("E22").Value = 4.4
$ws.Range("F22").Value = 281

# Row 23
$ws.Range("A23").Value = 32
$ws.Range("C23").Value = "Food Specialties Inc"
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0

# Row 24
$ws.Range("A24").Value = 40
$ws.Range("C24").Value = "Grand Park Sports Campus"
$ws.Range("E24").Value = 4.7
$ws.Range("F24").Value = 1377

# Row 25
$ws.Range("A25").Value = 49
$ws.Range("C25").Value = "Great Day Tattoo"
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = 18

# Row 26
$ws.Range("A26").Value = 58
$ws.Range("C26").Value = "HOA Affordable Fence Company"
$ws.Range("E26").Value = 4.1
$ws.Range("F26").Value = 106

# Row 27
$ws.Range("A27").Value = 16
$ws.Range("C27").Value = "Hamilton County Fairgrounds"
$ws.Range("F27").Value = 68

# Row 28
$ws.Range("A28").Value = 18
$ws.Range("C28").Value = "Harvest Pavillion"
$ws.Range("F28").Value = 18

# Row 29
$ws.Range("A29").Value = 31
$ws.Range("C29").Value = "Hendricks County Fairgrounds"
$ws.Range("E29").Value = 4.6
$ws.Range("F29").Value = 654

# Row 30
$ws.Range("A30").Value = 51
$ws.Range("C30").Value = "Hoosier Trim Products"
$ws.Range("E30").Value = 4.7
$ws.Range("F30").Value = 3

# Row 31
$ws.Range("A31").Value = 2
$ws.Range("C31").Value = "Indiana Black Expo Inc"
$ws.Range("E31").Value = 4.3
$ws.Range("F31").Value = 39

# Row 32
$ws.Range("A32").Value = 21
$ws.Range("C32").Value = "Indiana Convention Center"
$ws.Range("E32").Value = 4.5
$ws.Range("F32").Value = 528

# Row 33
$ws.Range("A33").Value = 20
$ws.Range("C33").Value = "Indiana Flower & Patio Show"
$ws.Range("E33").Value = 4.3
$ws.Range("F33").Value = 60

# Row 34
$ws.Range("A34").Value = 4
$ws.Range("C34").Value = "Indiana Latino Expo"
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0

# Row 35
$ws.Range("A35").Value = 15
$ws.Range("C35").Value = "Indiana State Fairgrounds & Event Center"
$ws.Range("E35").Value = 4.4
$ws.Range("F35").Value = 1344

# Row 36
$ws.Range("A36").Value = 36
$ws.Range("C36").Value = "Indiana State Numismatic Association"
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0

# Row 37
$ws.Range("A37").Value = 28
$ws.Range("C37").Value = "Indianapolis Auto Show"
$ws.Range("E37").Value = 3.1
$ws.Range("F37").Value = 51

# Row 38
$ws.Range("A38").Value = 9
$ws.Range("C38").Value = "Indianapolis Chapter of Indiana Black Expo, Inc."

# Row 39
$ws.Range("A39").Value = 23
$ws.Range("C39").Value = "Indianapolis Competition Products"
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 0

# Row 40
$ws.Range("A40").Value = 17
$ws.Range("C40").Value = "Indianapolis Motor Speedway"
$ws.Range("E40").Value = 4.8
$ws.Range("F40").Value = 11013

# Row 41
$ws.Range("A41").Value = 55
$ws.Range("C41").Value = "JW Marriott Indianapolis"
$ws.Range("E41").Value = 4.5
$ws.Range("F41").Value = 5755

# Row 42
$ws.Range("A42").Value = 19
$ws.Range("C42").Value = "Marketplace Events - Indianapolis Office"
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 1

# Row 43
$ws.Range("A43").Value = 6
$ws.Range("C43").Value = "Nail Expo"
$ws.Range("E43").Value = 3.2
$ws.Range("F43").Value = 62

# Row 44
$ws.Range("A44").Value = 8
$ws.Range("C44").Value = "Off Road Expo"
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0

# Row 45
$ws.Range("A45").Value = 50
$ws.Range("C45").Value = "Pan Am Tower"
$ws.Range("E45").Value = 4.3
$ws.Range("F45").Value = 28

# Row 46
$ws.Range("A46").Value = 33
$ws.Range("C46").Value = "Premier Surface"
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0

# Row 47
$ws.Range("A47").Value = 44
$ws.Range("C47").Value = "Purdue Extension / Horticulture Building"
$ws.Range("E47").Value = 4.4
$ws.Range("F47").Value = 36

# Row 48
$ws.Range("A48").Value = 52
$ws.Range("C48").Value = "Roberts Camera"
$ws.Range("E48").Value = 4.7
$ws.Range("F48").Value = 518

# Row 49
$ws.Range("A49").Value = 30
$ws.Range("C49").Value = "Royal Pin Western"
$ws.Range("E49").Value = 4.4
$ws.Range("F49").Value = 1182

# Row 50
$ws.Range("A50").Value = 39
$ws.Range("C50").Value = "Royal Pin Woodland"
$ws.Range("E50").Value = 4.3
$ws.Range("F50").Value = 1757

# Row 51
$ws.Range("A51").Value = 10
$ws.Range("C51").Value = "Samps Hack Shack Brownsburg"
$ws.Range("E51").Value = 4.9
$ws.Range("F51").Value = 22

# Row 52
$ws.Range("A52").Value = 13
$ws.Range("C52").Value = "Shepard Events"
$ws.Range("E52").Value = 3
$ws.Range("F52").Value = 1

# Row 53
$ws.Range("A53").Value = 34
$ws.Range("C53").Value = "Shepard Exposition Services"
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0

# Row 54
$ws.Range("A54").Value = 29
$ws.Range("C54").Value = "Suburban Indy Home & Outdoor Living Shows"
$ws.Range("E54").Value = 3.8
$ws.Range("F54").Value = 13

# Row 55
$ws.Range("A55").Value = 59
$ws.Range("C55").Value = "The Home Depot"
$ws.Range("E55").Value = 4.3
$ws.Range("F55").Value = 444

# Row 56
$ws.Range("A56").Value = 42
$ws.Range("C56").Value = "The Indiana Convention center"
$ws.Range("E56").Value = 4.7
$ws.Range("F56").Value = 21

# Row 57
$ws.Range("A57").Value = 47
$ws.Range("C57").Value = "The Korner Garage"
$ws.Range("F57").Value = 6
$ws.Range("D57").ClearContents()

# Row 58
$ws.Range("A58").Value = 54
$ws.Range("C58").Value = "Topgolf"
$ws.Range("E58").Value = 4.5
$ws.Range("F58").Value = 2966
$ws.Range("D58").Value = 2

